$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: finish the sentence ending in "message" with a period, then add
# the new "RTOS" paragraphs (Semaphores / Event groups / Message queues),
# a blank paragraph, and the new "Retarget IO" section - all inserted right
# before the existing (empty) paragraph that precedes "Ez I2C".
# ---------------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("Free RTSO Task was created with semaphores to run UART task on key press. And print a message", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter(".")

$r.Collapse(0)
$r.InsertAfter("`rSemaphores are enabled and use for ease of resource sharing between tasks.")

$r.Collapse(0)
$r.InsertAfter("`rEvent groups are used to notify observation tasks like ezI2c and BLE of the change in the motor PWM duty cycle which can then be communicated to user.")

$r.Collapse(0)
$r.InsertAfter("`rMessage queues are used to communicate the pwm duty cycle percentage between tasks, mainly to communicate the duty cycle to be set from capsense task to the motor control task.")

$r.Collapse(0)
$r.InsertAfter("`r")

$r.Collapse(0)
$r.InsertAfter("`rRetarget IO")

$r.Collapse(0)
$r.InsertAfter("`rThe peripheral driver library " + [char]0x201C + "Retarget io" + [char]0x201D + " was enabled in the projects build settings to re route the messages from printf() to UART on SCB5 for displaying the messages in terminal application.")

# ---------------------------------------------------------------------------
# Change 2: split the last paragraph right after "...0 to 100." so the
# trailing bookmark / line-break run move into a new paragraph that starts
# with the new sentence about the motor control task message queue.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$null = $r2.Find.Execute("duty cycle for the motor in percent from 0 to 100.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Collapse(0)

$d.Bookmarks("_GoBack").Delete()

$r2.InsertParagraphAfter()
$r2.Collapse(0)
$r2.MoveStart(1, 1)
$r2.InsertAfter("The duty cycle is then communicated to the motor control task via message queue for the same to set the duty cycle of the motor to this perncentage.")
$r2.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r2)

Write-Host "done"
